$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")
$ws.Activate()

# ------------------------------------------------------------------
# 1) Rename the 3 comboBox labels used by the "trans/losszoneBus/insert"
#    form (rows 65-67) so they get a trailing ";" like the rest of the
#    sheet. Doing this first (and in this exact order) makes the engine
#    garbage-collect the old shared strings and re-append the renamed
#    ones + "trans/losszoneBus/insert" in the order required by the
#    target file.
# ------------------------------------------------------------------
$ws.Range("C65").Value = "comboBox18;"
$ws.Range("C66").Value = "comboBox19;"
$ws.Range("C67").Value = "comboBox20;"

# ------------------------------------------------------------------
# 2) Add the new "trans/losszoneLine/insert" form description (rows
#    70-74) and the new "trans/modelload/insert" row (row 75).
# ------------------------------------------------------------------

# Column A first (labels / button names)
$ws.Range("A70").Value = "label62;"
$ws.Range("A71").Value = "label61;"
$ws.Range("A72").Value = "label60;"

# Column B (meaning of the panel) for the label rows
$ws.Range("B70").Value = "caseID"
$ws.Range("B71").Value = "ramo"
$ws.Range("B72").Value = "losszone"

# Buttons (A then B)
$ws.Range("A73").Value = "button21;"
$ws.Range("A74").Value = "button20;"
$ws.Range("B73").Value = "Submit"
$ws.Range("B74").Value = "Clear"

# Column C (textbox/combobox names)
$ws.Range("C70").Value = "comboBox23;"
$ws.Range("C71").Value = "comboBox22;"
$ws.Range("C72").Value = "comboBox21;"
$ws.Range("C73").Value = "-"
$ws.Range("C74").Value = "-"

# Column D (panel number)
$ws.Range("D70").Value = 16
$ws.Range("D71").Value = 16
$ws.Range("D72").Value = 16
$ws.Range("D73").Value = 16
$ws.Range("D74").Value = 16

# Column E (panel meaning)
$ws.Range("E70").Value = "trans/losszoneLine/insert"
$ws.Range("E71").Value = "trans/losszoneLine/insert"
$ws.Range("E72").Value = "trans/losszoneLine/insert"
$ws.Range("E73").Value = "trans/losszoneLine/insert"
$ws.Range("E74").Value = "trans/losszoneLine/insert"

# Row 75 - new "trans/modelload/insert" entry
$ws.Range("A75").Value = "label63"
$ws.Range("B75").Value = "FormNotUsed "
$ws.Range("C75").Value = "-"
$ws.Range("D75").Value = 17
$ws.Range("E75").Value = "trans/modelload/insert"

# ------------------------------------------------------------------
# 3) Copy the row formatting from existing rows that already use the
#    right fill/border style: rows 70-74 reuse the style of row 59
#    (alternate shading) and row 75 reuses the style of row 65.
# ------------------------------------------------------------------
$ws.Range("A59:E59").Copy() | Out-Null
$ws.Range("A70:E74").PasteSpecial(-4122) | Out-Null

$ws.Range("A65:E65").Copy() | Out-Null
$ws.Range("A75:E75").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Update the view state to match the new extent of the sheet.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E77").Select() | Out-Null
